# Apply FlashScore odds updates for 2024-10-21 workbook (commit: "Atualizando o arquivo XLSX")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.67
$ws.Range("U2").Value = 1.95
$ws.Range("V2").Value = 1.8
$ws.Range("BC2").Value = 130
$ws.Range("U3").Value = 1.95
$ws.Range("V3").Value = 1.8
$ws.Range("V4").Value = 1.67
$ws.Range("Q5").Value = 2.6
$ws.Range("R5").Value = 1.48
$ws.Range("V5").Value = 1.67
$ws.Range("G11").Value = 1.7
$ws.Range("I11").Value = 4.33
$ws.Range("N11").Value = 17
$ws.Range("O11").Value = 1.17
$ws.Range("P11").Value = 5
$ws.Range("Q11").Value = 1.57
$ws.Range("R11").Value = 2.35
$ws.Range("U11").Value = 1.57
$ws.Range("V11").Value = 2.25
$ws.Range("Z11").Value = 15
$ws.Range("AC11").Value = 17
$ws.Range("AG11").Value = 126
$ws.Range("AX11").Value = 21
$ws.Range("G12").Value = 2.15
$ws.Range("I12").Value = 3.25
$ws.Range("L12").Value = 3.6
$ws.Range("Q12").Value = 1.67
$ws.Range("R12").Value = 2.15
$ws.Range("AA12").Value = 15
$ws.Range("AB12").Value = 21
$ws.Range("AI12").Value = 19
$ws.Range("AN12").Value = 4.33
$ws.Range("AO12").Value = 11
$ws.Range("AW12").Value = 5.5
$ws.Range("AX12").Value = 17
$ws.Range("BA12").Value = 67
$ws.Range("G16").Value = 1.5
$ws.Range("H16").Value = 4.33
$ws.Range("I16").Value = 5.75
$ws.Range("J16").Value = 2.05
$ws.Range("N16").Value = 12
$ws.Range("S16").Value = 1.36
$ws.Range("T16").Value = 3
$ws.Range("U16").Value = 1.91
$ws.Range("V16").Value = 1.8
$ws.Range("X16").Value = 7
$ws.Range("AE16").Value = 19
$ws.Range("AK16").Value = 67
$ws.Range("AT16").Value = 3
$ws.Range("AX16").Value = 34
$ws.Range("G17").Value = 2.5
$ws.Range("H17").Value = 3.8
$ws.Range("I17").Value = 2.45
$ws.Range("J17").Value = 3
$ws.Range("L17").Value = 2.88
$ws.Range("N17").Value = 21
$ws.Range("S17").Value = 1.22
$ws.Range("T17").Value = 4
$ws.Range("AD17").Value = 8.5
$ws.Range("AM17").Value = 19
$ws.Range("AT17").Value = 4
$ws.Range("AV17").Value = 34
$ws.Range("AX17").Value = 12
$ws.Range("I20").Value = 1.41
$ws.Range("G21").Value = 1.33
$ws.Range("M23").Value = 1.08
$ws.Range("O23").Value = 1.44
$ws.Range("P23").Value = 2.63
$ws.Range("M24").Value = 1.08
$ws.Range("O24").Value = 1.4
$ws.Range("M25").Value = 1.05
$ws.Range("O25").Value = 1.25
$ws.Range("M26").Value = 1.04
$ws.Range("O26").Value = 1.22
$ws.Range("Q26").Value = 1.77
$ws.Range("M30").Value = 1.07
$ws.Range("O30").Value = 1.3
$ws.Range("G32").Value = 2.25
$ws.Range("I32").Value = 2.88
$ws.Range("J32").Value = 2.87
$ws.Range("K32").Value = 2.37
$ws.Range("L32").Value = 3.25
$ws.Range("M32").Value = 1.03
$ws.Range("O32").Value = 1.17
$ws.Range("AY32").Value = 21
$ws.Range("I33").Value = 2.3
$ws.Range("M33").Value = 1.08
$ws.Range("O33").Value = 1.4
$ws.Range("G34").Value = 1.6
$ws.Range("U34").Value = 1.91
$ws.Range("V34").Value = 1.91
$ws.Range("X34").Value = 7.5
$ws.Range("AO34").Value = 8
$ws.Range("AU34").Value = 8.5
$ws.Range("AX34").Value = 29
$ws.Range("G35").Value = 2.4
$ws.Range("I35").Value = 3.2
$ws.Range("J35").Value = 3.1
$ws.Range("L35").Value = 3.75
$ws.Range("U35").Value = 1.83
$ws.Range("V35").Value = 1.83
$ws.Range("AO35").Value = 13
$ws.Range("AQ35").Value = 41
$ws.Range("AW35").Value = 5
$ws.Range("U36").Value = 1.73
$ws.Range("V38").Value = 1.73
$ws.Range("U39").Value = 1.73
